$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = '67.356.61'
$ws.Cells.Item(2, 5).Value2 = '  -0.14%  '
$ws.Cells.Item(3, 4).Value2 = '3.482.45'
$ws.Cells.Item(3, 5).Value2 = '  -1.04%  '
$ws.Cells.Item(4, 5).Value2 = '  -0.01%  '
$ws.Cells.Item(5, 4).Value2 = "'604.23"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value2 = '  -1.83%  '
$ws.Cells.Item(6, 4).Value2 = "'150.85"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value2 = '  -0.43%  '
$ws.Cells.Item(7, 4).Value2 = '3.478.58'
$ws.Cells.Item(7, 5).Value2 = '  -1.08%  '
$ws.Cells.Item(8, 5).Value2 = '  +0.05%  '
$ws.Cells.Item(9, 4).Value2 = "'0.485"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value2 = '  +0.86%  '
$ws.Cells.Item(10, 5).Value2 = '  +2.52%  '
$ws.Cells.Item(11, 4).Value2 = "'7.54"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value2 = '  +5.82%  '
$ws.Cells.Item(12, 4).Value2 = "'0.430"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value2 = '  +0.93%  '
$ws.Cells.Item(13, 5).Value2 = '  -2.17%  '
$ws.Cells.Item(14, 4).Value2 = "'32.01"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value2 = '  -0.36%  '
$ws.Cells.Item(15, 4).Value2 = '4.073.40'
$ws.Cells.Item(15, 5).Value2 = '  -1.01%  '
$ws.Cells.Item(16, 4).Value2 = '3.485.68'
$ws.Cells.Item(16, 5).Value2 = '  -0.96%  '
$ws.Cells.Item(17, 4).Value2 = '67.267.25'
$ws.Cells.Item(17, 5).Value2 = '  -0.26%  '
$ws.Cells.Item(18, 5).Value2 = '  -0.84%  '
$ws.Cells.Item(19, 5).Value2 = '  +1.19%  '
$ws.Cells.Item(20, 4).Value2 = "'15.30"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value2 = '  -0.44%  '
$ws.Cells.Item(21, 5).Value2 = '  +3.15%  '
$ws.Cells.Item(22, 4).Value2 = "'445.08"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value2 = '  +0.16%  '
$ws.Cells.Item(23, 4).Value2 = "'0.626"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value2 = '  +0.42%  '
$ws.Cells.Item(24, 4).Value2 = "'78.12"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value2 = '  +0.96%  '
$ws.Cells.Item(25, 4).Value2 = '3.628.41'
$ws.Cells.Item(25, 5).Value2 = '  -0.86%  '
$ws.Cells.Item(26, 5).Value2 = '  -0.05%  '
$ws.Cells.Item(27, 5).Value2 = '  -3.41%  '
$ws.Cells.Item(28, 4).Value2 = "'8.72"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value2 = '  +2.45%  '
$ws.Cells.Item(29, 4).Value2 = "'9.96"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value2 = '  -2.63%  '
$ws.Cells.Item(30, 4).Value2 = "'2.49"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value2 = '  -0.97%  '
$ws.Cells.Item(31, 5).Value2 = '  +2.99%  '
$ws.Cells.Item(32, 4).Value2 = "'0.171"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value2 = '  +4.07%  '
$ws.Cells.Item(33, 5).Value2 = '  -0.07%  '
$ws.Cells.Item(34, 4).Value2 = "'25.48"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value2 = '  -1.44%  '
$ws.Cells.Item(35, 4).Value2 = "'6.11"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value2 = '  -0.68%  '
$ws.Cells.Item(36, 5).Value2 = '  +0.38%  '
$ws.Cells.Item(37, 4).Value2 = '3.478.05'
$ws.Cells.Item(37, 5).Value2 = '  -0.99%  '
$ws.Cells.Item(38, 4).Value2 = "'7.95"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value2 = '  -0.80%  '
$ws.Cells.Item(39, 5).Value2 = '  +0.01%  '
$ws.Cells.Item(40, 5).Value2 = '  +7.00%  '
$ws.Cells.Item(41, 2).Value2 = 'FirstDigitalUSD'
$ws.Cells.Item(41, 3).Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(41, 4).Value2 = "'1.00"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value2 = '  +0.01%  '
$ws.Cells.Item(42, 2).Value2 = 'Monero'
$ws.Cells.Item(42, 3).Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(42, 4).Value2 = "'177.31"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value2 = '  -0.09%  '
$ws.Cells.Item(43, 4).Value2 = "'0.0891"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value2 = '  +0.91%  '
$ws.Cells.Item(44, 4).Value2 = "'5.41"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value2 = '  -0.09%  '
$ws.Cells.Item(45, 4).Value2 = "'0.888"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value2 = '  +0.78%  '
$ws.Cells.Item(46, 4).Value2 = "'29.97"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value2 = '  +5.73%  '
$ws.Cells.Item(47, 4).Value2 = "'46.38"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value2 = '  +3.06%  '
$ws.Cells.Item(48, 5).Value2 = '  +3.37%  '
$ws.Cells.Item(49, 4).Value2 = "'2.53"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value2 = '  -3.66%  '
$ws.Cells.Item(50, 4).Value2 = "'7.58"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value2 = '  -0.37%  '
$ws.Cells.Item(51, 4).Value2 = "'0.983"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value2 = '  -1.12%  '
